$d = $word.ActiveDocument

# 1. "LChakra levels" -> "Chakra levels"
$d.Content.Find.Execute("LChakra levels", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Chakra levels", 2)

# 2. "How much chakra a ninja has to expend." -> "How much chakra a ninja has available to expend."
$d.Content.Find.Execute("How much chakra a ninja has to expend.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "How much chakra a ninja has available to expend.", 2)
